# Debug | Fix save data to file | Split traj to xyz in trials table
#
# The trials table used to store each body-trajectory as a single combined
# column ("<name>_traj" / "trajectory (X,Y,Z) in time"). This change splits
# each of those three trajectory columns (target_traj, prime_traj, pas_traj)
# into three separate columns: "<name>_x", "<name>_y", "<name>_z"
# (each documented as "X/Y/Z trajectory in time"). Everything to the right
# shifts over to make room, which Insert() handles for us (carrying styles
# and the existing cat_block column-width setting along automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- split target_traj (originally column AA) into target_x / target_y / target_z ---
$ws.Columns("AA:AB").Insert()
$ws.Range("AA1").Value = "target_x"
$ws.Range("AB1").Value = "target_y"
$ws.Range("AC1").Value = "target_z"
$ws.Range("AA2").Value = "X trajectory in time"
$ws.Range("AB2").Value = "Y trajectory in time"
$ws.Range("AC2").Value = "Z trajectory in time"

# --- split prime_traj (originally column AG, now shifted to AI) into prime_x / prime_y / prime_z ---
$ws.Columns("AI:AJ").Insert()
$ws.Range("AI1").Value = "prime_x"
$ws.Range("AJ1").Value = "prime_y"
$ws.Range("AK1").Value = "prime_z"
$ws.Range("AI2").Value = "X trajectory in time"
$ws.Range("AJ2").Value = "Y trajectory in time"
$ws.Range("AK2").Value = "Z trajectory in time"

# --- split pas_traj (originally column AL, now shifted to AP) into pas_x / pas_y / pas_z ---
$ws.Columns("AP:AQ").Insert()
$ws.Range("AP1").Value = "pas_x"
$ws.Range("AQ1").Value = "pas_y"
$ws.Range("AR1").Value = "pas_z"
$ws.Range("AP2").Value = "X trajectory in time"
$ws.Range("AQ2").Value = "Y trajectory in time"
$ws.Range("AR2").Value = "Z trajectory in time"

# update view state: selection moves to AC2, and the visible window scrolls
# so column Y is the left-most visible column
$ws.Range("AC2").Select()
$excel.ActiveWindow.ScrollColumn = 25
$excel.ActiveWindow.ScrollRow = 1

Write-Output "Split trajectory columns into X/Y/Z for target, prime and pas."
